# Generate Report for Archive
# - Update localization status text "Ready for handoff" -> "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns E & F) and on the
#   per-locale sheets (zh-cn, de-de) "Status" column (C).
# - Narrow the now-shorter "Status" columns to better fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column C ---
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

$dede.Columns.Item(3).ColumnWidth = 12.5
